$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44466
$ws.Cells.Item(2, 10).Value = 1150
$ws.Cells.Item(2, 11).Value = 14000
$ws.Cells.Item(2, 12).Value = 15000
$ws.Cells.Item(2, 13).Value = 14500
$ws.Cells.Item(2, 16).Value = 1115
$ws.Cells.Item(3, 4).Value = 44466
$ws.Cells.Item(3, 9).Value = "Segunda"
$ws.Cells.Item(3, 10).Value = 790
$ws.Cells.Item(3, 11).Value = 12000
$ws.Cells.Item(3, 12).Value = 12000
$ws.Cells.Item(3, 13).Value = 12000
$ws.Cells.Item(3, 16).Value = 923
$ws.Cells.Item(4, 4).Value = 44445
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 790
$ws.Cells.Item(4, 11).Value = 13000
$ws.Cells.Item(4, 12).Value = 14000
$ws.Cells.Item(4, 13).Value = 13494
$ws.Cells.Item(4, 16).Value = 1038
$ws.Cells.Item(5, 4).Value = 44445
$ws.Cells.Item(5, 9).Value = "Segunda"
$ws.Cells.Item(5, 10).Value = 340
$ws.Cells.Item(5, 11).Value = 11000
$ws.Cells.Item(5, 12).Value = 12000
$ws.Cells.Item(5, 13).Value = 11500
$ws.Cells.Item(5, 16).Value = 885
$ws.Cells.Item(6, 4).Value = 44571
$ws.Cells.Item(7, 4).Value = 44571
$ws.Cells.Item(7, 9).Value = "Segunda"
$ws.Cells.Item(7, 10).Value = 106
$ws.Cells.Item(7, 11).Value = 10000
$ws.Cells.Item(7, 12).Value = 10000
$ws.Cells.Item(7, 13).Value = 10000
$ws.Cells.Item(7, 16).Value = 769
$ws.Cells.Item(8, 4).Value = 44165
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 430
$ws.Cells.Item(8, 11).Value = 31000
$ws.Cells.Item(8, 12).Value = 32000
$ws.Cells.Item(8, 13).Value = 31465
$ws.Cells.Item(8, 16).Value = 2420
$ws.Cells.Item(9, 4).Value = 44606
$ws.Cells.Item(9, 10).Value = 520
$ws.Cells.Item(9, 11).Value = 17000
$ws.Cells.Item(9, 12).Value = 18000
$ws.Cells.Item(9, 13).Value = 17500
$ws.Cells.Item(9, 16).Value = 1346
$ws.Cells.Item(10, 4).Value = 44263
$ws.Cells.Item(10, 10).Value = 250
$ws.Cells.Item(10, 11).Value = 40000
$ws.Cells.Item(10, 12).Value = 40000
$ws.Cells.Item(10, 13).Value = 40000
$ws.Cells.Item(10, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(10, 16).Value = 2667
$ws.Cells.Item(10, 17).Value = 15
$ws.Cells.Item(11, 4).Value = 44515
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 1060
$ws.Cells.Item(11, 11).Value = 16000
$ws.Cells.Item(11, 12).Value = 18000
$ws.Cells.Item(11, 13).Value = 17000
$ws.Cells.Item(11, 16).Value = 1308
$ws.Cells.Item(12, 4).Value = 44515
$ws.Cells.Item(12, 9).Value = "Segunda"
$ws.Cells.Item(12, 11).Value = 14000
$ws.Cells.Item(12, 12).Value = 14000
$ws.Cells.Item(12, 13).Value = 14000
$ws.Cells.Item(12, 16).Value = 1077
$ws.Cells.Item(13, 4).Value = 44435
$ws.Cells.Item(13, 10).Value = 880
$ws.Cells.Item(13, 11).Value = 13000
$ws.Cells.Item(13, 12).Value = 14000
$ws.Cells.Item(13, 13).Value = 13500
$ws.Cells.Item(13, 16).Value = 1038
$ws.Cells.Item(14, 4).Value = 44435
$ws.Cells.Item(14, 11).Value = 11000
$ws.Cells.Item(14, 12).Value = 12000
$ws.Cells.Item(14, 13).Value = 11500
$ws.Cells.Item(14, 16).Value = 885
$ws.Cells.Item(15, 4).Value = 44270
$ws.Cells.Item(16, 4).Value = 44536
$ws.Cells.Item(16, 10).Value = 790
$ws.Cells.Item(16, 11).Value = 14000
$ws.Cells.Item(16, 12).Value = 15000
$ws.Cells.Item(16, 13).Value = 14494
$ws.Cells.Item(16, 16).Value = 1115
$ws.Cells.Item(17, 4).Value = 44536
$ws.Cells.Item(17, 10).Value = 430
$ws.Cells.Item(17, 11).Value = 11000
$ws.Cells.Item(17, 12).Value = 12000
$ws.Cells.Item(17, 13).Value = 11500
$ws.Cells.Item(17, 16).Value = 885
$ws.Cells.Item(18, 4).Value = 44354
$ws.Cells.Item(18, 10).Value = 250
$ws.Cells.Item(18, 11).Value = 15000
$ws.Cells.Item(18, 12).Value = 16000
$ws.Cells.Item(18, 13).Value = 15500
$ws.Cells.Item(18, 16).Value = 1192
$ws.Cells.Item(19, 4).Value = 44662
$ws.Cells.Item(19, 10).Value = 610
$ws.Cells.Item(19, 11).Value = 12000
$ws.Cells.Item(19, 12).Value = 13000
$ws.Cells.Item(19, 13).Value = 12500
$ws.Cells.Item(19, 16).Value = 962
$ws.Cells.Item(20, 4).Value = 44438
$ws.Cells.Item(20, 11).Value = 13000
$ws.Cells.Item(20, 12).Value = 14000
$ws.Cells.Item(20, 13).Value = 13494
$ws.Cells.Item(20, 16).Value = 1038
$ws.Cells.Item(21, 4).Value = 44438
$ws.Cells.Item(21, 9).Value = "Segunda"
$ws.Cells.Item(21, 10).Value = 340
$ws.Cells.Item(21, 11).Value = 11000
$ws.Cells.Item(21, 12).Value = 12000
$ws.Cells.Item(21, 13).Value = 11500
$ws.Cells.Item(21, 16).Value = 885
$ws.Cells.Item(22, 4).Value = 44613
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 790
$ws.Cells.Item(22, 11).Value = 16000
$ws.Cells.Item(22, 12).Value = 17000
$ws.Cells.Item(22, 13).Value = 16494
$ws.Cells.Item(22, 16).Value = 1269
$ws.Cells.Item(23, 4).Value = 44179
$ws.Cells.Item(23, 10).Value = 430
$ws.Cells.Item(23, 11).Value = 29000
$ws.Cells.Item(23, 12).Value = 30000
$ws.Cells.Item(23, 13).Value = 29465
$ws.Cells.Item(23, 16).Value = 2267
$ws.Cells.Item(24, 4).Value = 44371
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 160
$ws.Cells.Item(24, 11).Value = 20000
$ws.Cells.Item(24, 12).Value = 21000
$ws.Cells.Item(24, 13).Value = 20500
$ws.Cells.Item(24, 16).Value = 1577
$ws.Cells.Item(25, 4).Value = 44655
$ws.Cells.Item(25, 11).Value = 14000
$ws.Cells.Item(25, 12).Value = 15000
$ws.Cells.Item(25, 13).Value = 14494
$ws.Cells.Item(25, 16).Value = 1115
$ws.Cells.Item(26, 4).Value = 44410
$ws.Cells.Item(26, 10).Value = 790
$ws.Cells.Item(26, 11).Value = 15000
$ws.Cells.Item(26, 12).Value = 16000
$ws.Cells.Item(26, 13).Value = 15494
$ws.Cells.Item(26, 16).Value = 1192
$ws.Cells.Item(27, 4).Value = 44410
$ws.Cells.Item(27, 10).Value = 340
$ws.Cells.Item(27, 11).Value = 13000
$ws.Cells.Item(27, 12).Value = 13000
$ws.Cells.Item(27, 13).Value = 13000
$ws.Cells.Item(27, 16).Value = 1000
$ws.Cells.Item(28, 4).Value = 44417
$ws.Cells.Item(28, 10).Value = 790
$ws.Cells.Item(28, 11).Value = 14000
$ws.Cells.Item(28, 12).Value = 15000
$ws.Cells.Item(28, 13).Value = 14500
$ws.Cells.Item(28, 16).Value = 1115
$ws.Cells.Item(29, 4).Value = 44417
$ws.Cells.Item(29, 9).Value = "Segunda"
$ws.Cells.Item(29, 10).Value = 340
$ws.Cells.Item(29, 11).Value = 13000
$ws.Cells.Item(29, 12).Value = 13000
$ws.Cells.Item(29, 13).Value = 13000
$ws.Cells.Item(29, 16).Value = 1000
$ws.Cells.Item(30, 4).Value = 44676
$ws.Cells.Item(30, 10).Value = 790
$ws.Cells.Item(30, 11).Value = 11000
$ws.Cells.Item(30, 12).Value = 12000
$ws.Cells.Item(30, 13).Value = 11494
$ws.Cells.Item(30, 16).Value = 884
$ws.Cells.Item(31, 4).Value = 44522
$ws.Cells.Item(31, 11).Value = 16000
$ws.Cells.Item(31, 12).Value = 18000
$ws.Cells.Item(31, 13).Value = 16987
$ws.Cells.Item(31, 16).Value = 1307
$ws.Cells.Item(32, 4).Value = 44522
$ws.Cells.Item(32, 9).Value = "Segunda"
$ws.Cells.Item(32, 10).Value = 360
$ws.Cells.Item(32, 11).Value = 15000
$ws.Cells.Item(32, 12).Value = 15000
$ws.Cells.Item(32, 13).Value = 15000
$ws.Cells.Item(32, 16).Value = 1154
$ws.Cells.Item(33, 4).Value = 44333
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 11).Value = 25000
$ws.Cells.Item(33, 12).Value = 26000
$ws.Cells.Item(33, 13).Value = 25500
$ws.Cells.Item(33, 16).Value = 1962
$ws.Cells.Item(34, 4).Value = 44333
$ws.Cells.Item(34, 9).Value = "Segunda"
$ws.Cells.Item(34, 10).Value = 160
$ws.Cells.Item(34, 11).Value = 23000
$ws.Cells.Item(34, 12).Value = 23000
$ws.Cells.Item(34, 13).Value = 23000
$ws.Cells.Item(34, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(34, 16).Value = 1769
$ws.Cells.Item(34, 17).Value = 13
$ws.Cells.Item(35, 4).Value = 44277
$ws.Cells.Item(35, 11).Value = 38000
$ws.Cells.Item(35, 12).Value = 38000
$ws.Cells.Item(35, 13).Value = 38000
$ws.Cells.Item(35, 16).Value = 2923
$ws.Cells.Item(36, 4).Value = 44277
$ws.Cells.Item(36, 11).Value = 35000
$ws.Cells.Item(36, 12).Value = 35000
$ws.Cells.Item(36, 13).Value = 35000
$ws.Cells.Item(36, 16).Value = 2692
$ws.Cells.Item(37, 4).Value = 44186
$ws.Cells.Item(37, 10).Value = 450
$ws.Cells.Item(37, 11).Value = 29000
$ws.Cells.Item(37, 12).Value = 30000
$ws.Cells.Item(37, 13).Value = 29444
$ws.Cells.Item(37, 16).Value = 2265
$ws.Cells.Item(38, 4).Value = 44326
$ws.Cells.Item(38, 10).Value = 340
$ws.Cells.Item(38, 11).Value = 25000
$ws.Cells.Item(38, 12).Value = 25000
$ws.Cells.Item(38, 13).Value = 25000
$ws.Cells.Item(38, 16).Value = 1923
$ws.Cells.Item(39, 4).Value = 44326
$ws.Cells.Item(39, 11).Value = 23000
$ws.Cells.Item(39, 12).Value = 23000
$ws.Cells.Item(39, 13).Value = 23000
$ws.Cells.Item(39, 16).Value = 1769
$ws.Cells.Item(40, 4).Value = 44340
$ws.Cells.Item(40, 10).Value = 250
$ws.Cells.Item(40, 11).Value = 20000
$ws.Cells.Item(40, 12).Value = 20000
$ws.Cells.Item(40, 13).Value = 20000
$ws.Cells.Item(40, 16).Value = 1538
$ws.Cells.Item(41, 4).Value = 44340
$ws.Cells.Item(41, 11).Value = 18000
$ws.Cells.Item(41, 12).Value = 18000
$ws.Cells.Item(41, 13).Value = 18000
$ws.Cells.Item(41, 16).Value = 1385
$ws.Cells.Item(42, 4).Value = 44585
$ws.Cells.Item(42, 10).Value = 790
$ws.Cells.Item(42, 11).Value = 10000
$ws.Cells.Item(42, 12).Value = 11000
$ws.Cells.Item(42, 13).Value = 10494
$ws.Cells.Item(42, 16).Value = 807
$ws.Cells.Item(43, 4).Value = 44585
$ws.Cells.Item(43, 9).Value = "Segunda"
$ws.Cells.Item(43, 11).Value = 9000
$ws.Cells.Item(43, 12).Value = 9000
$ws.Cells.Item(43, 13).Value = 9000
$ws.Cells.Item(43, 16).Value = 692
$ws.Cells.Item(44, 4).Value = 44592
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 400
$ws.Cells.Item(44, 11).Value = 9000
$ws.Cells.Item(44, 12).Value = 10000
$ws.Cells.Item(44, 13).Value = 9575
$ws.Cells.Item(44, 16).Value = 737
$ws.Cells.Item(45, 4).Value = 44592
$ws.Cells.Item(45, 9).Value = "Segunda"
$ws.Cells.Item(45, 10).Value = 100
$ws.Cells.Item(45, 11).Value = 8000
$ws.Cells.Item(45, 12).Value = 8000
$ws.Cells.Item(45, 13).Value = 8000
$ws.Cells.Item(45, 16).Value = 615
$ws.Cells.Item(46, 4).Value = 44648
$ws.Cells.Item(46, 10).Value = 610
$ws.Cells.Item(46, 11).Value = 16000
$ws.Cells.Item(46, 12).Value = 17000
$ws.Cells.Item(46, 13).Value = 16500
$ws.Cells.Item(46, 16).Value = 1269
$ws.Cells.Item(47, 4).Value = 44242
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 250
$ws.Cells.Item(47, 11).Value = 41000
$ws.Cells.Item(47, 12).Value = 43000
$ws.Cells.Item(47, 13).Value = 42000
$ws.Cells.Item(47, 16).Value = 3231
$ws.Cells.Item(48, 4).Value = 44620
$ws.Cells.Item(48, 10).Value = 790
$ws.Cells.Item(48, 11).Value = 15000
$ws.Cells.Item(48, 12).Value = 16000
$ws.Cells.Item(48, 13).Value = 15494
$ws.Cells.Item(48, 16).Value = 1192
$ws.Cells.Item(49, 4).Value = 44627
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 790
$ws.Cells.Item(49, 11).Value = 14000
$ws.Cells.Item(49, 12).Value = 15000
$ws.Cells.Item(49, 13).Value = 14494
$ws.Cells.Item(49, 16).Value = 1115
$ws.Cells.Item(50, 4).Value = 44627
$ws.Cells.Item(50, 9).Value = "Segunda"
$ws.Cells.Item(50, 10).Value = 340
$ws.Cells.Item(50, 11).Value = 13000
$ws.Cells.Item(50, 12).Value = 13000
$ws.Cells.Item(50, 13).Value = 13000
$ws.Cells.Item(50, 16).Value = 1000
$ws.Cells.Item(51, 4).Value = 44298
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 340
$ws.Cells.Item(51, 11).Value = 24000
$ws.Cells.Item(51, 12).Value = 25000
$ws.Cells.Item(51, 13).Value = 24500
$ws.Cells.Item(51, 16).Value = 1885
$ws.Cells.Item(52, 4).Value = 44459
$ws.Cells.Item(52, 10).Value = 970
$ws.Cells.Item(52, 11).Value = 13000
$ws.Cells.Item(52, 12).Value = 14000
$ws.Cells.Item(52, 13).Value = 13495
$ws.Cells.Item(52, 16).Value = 1038
$ws.Cells.Item(53, 4).Value = 44459
$ws.Cells.Item(53, 9).Value = "Segunda"
$ws.Cells.Item(53, 10).Value = 520
$ws.Cells.Item(53, 11).Value = 11000
$ws.Cells.Item(53, 12).Value = 12000
$ws.Cells.Item(53, 13).Value = 11500
$ws.Cells.Item(53, 16).Value = 885
$ws.Cells.Item(54, 4).Value = 44312
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 430
$ws.Cells.Item(54, 11).Value = 25000
$ws.Cells.Item(54, 12).Value = 25000
$ws.Cells.Item(54, 13).Value = 25000
$ws.Cells.Item(54, 16).Value = 1923
$ws.Cells.Item(55, 4).Value = 44312
$ws.Cells.Item(55, 9).Value = "Segunda"
$ws.Cells.Item(55, 10).Value = 250
$ws.Cells.Item(55, 11).Value = 23000
$ws.Cells.Item(55, 12).Value = 23000
$ws.Cells.Item(55, 13).Value = 23000
$ws.Cells.Item(55, 16).Value = 1769
$ws.Cells.Item(56, 4).Value = 44172
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 430
$ws.Cells.Item(56, 11).Value = 30000
$ws.Cells.Item(56, 12).Value = 30000
$ws.Cells.Item(56, 13).Value = 30000
$ws.Cells.Item(56, 16).Value = 2308
$ws.Cells.Item(57, 4).Value = 44319
$ws.Cells.Item(57, 10).Value = 340
$ws.Cells.Item(57, 11).Value = 24000
$ws.Cells.Item(57, 12).Value = 25000
$ws.Cells.Item(57, 13).Value = 24500
$ws.Cells.Item(57, 16).Value = 1885
$ws.Cells.Item(58, 4).Value = 44319
$ws.Cells.Item(58, 10).Value = 160
$ws.Cells.Item(58, 11).Value = 22000
$ws.Cells.Item(58, 12).Value = 22000
$ws.Cells.Item(58, 13).Value = 22000
$ws.Cells.Item(58, 16).Value = 1692
$ws.Cells.Item(59, 4).Value = 44508
$ws.Cells.Item(59, 10).Value = 1150
$ws.Cells.Item(59, 11).Value = 15000
$ws.Cells.Item(59, 12).Value = 16000
$ws.Cells.Item(59, 13).Value = 15500
$ws.Cells.Item(59, 16).Value = 1192
$ws.Cells.Item(60, 4).Value = 44291
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 340
$ws.Cells.Item(60, 11).Value = 24000
$ws.Cells.Item(60, 12).Value = 25000
$ws.Cells.Item(60, 13).Value = 24500
$ws.Cells.Item(60, 16).Value = 1885
$ws.Cells.Item(61, 4).Value = 44690
$ws.Cells.Item(61, 10).Value = 790
$ws.Cells.Item(61, 11).Value = 12000
$ws.Cells.Item(61, 12).Value = 13000
$ws.Cells.Item(61, 13).Value = 12494
$ws.Cells.Item(61, 16).Value = 961
$ws.Cells.Item(62, 4).Value = 44641
$ws.Cells.Item(62, 10).Value = 610
$ws.Cells.Item(62, 11).Value = 14000
$ws.Cells.Item(62, 12).Value = 15000
$ws.Cells.Item(62, 13).Value = 14500
$ws.Cells.Item(62, 16).Value = 1115
$ws.Cells.Item(63, 4).Value = 44221
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 350
$ws.Cells.Item(63, 11).Value = 40000
$ws.Cells.Item(63, 12).Value = 42000
$ws.Cells.Item(63, 13).Value = 40857
$ws.Cells.Item(63, 16).Value = 3143
$ws.Cells.Item(64, 4).Value = 44221
$ws.Cells.Item(64, 9).Value = "Segunda"
$ws.Cells.Item(64, 10).Value = 180
$ws.Cells.Item(64, 11).Value = 35000
$ws.Cells.Item(64, 12).Value = 35000
$ws.Cells.Item(64, 13).Value = 35000
$ws.Cells.Item(64, 16).Value = 2692
$ws.Cells.Item(65, 4).Value = 44494
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 780
$ws.Cells.Item(65, 11).Value = 15000
$ws.Cells.Item(65, 12).Value = 15000
$ws.Cells.Item(65, 13).Value = 15000
$ws.Cells.Item(65, 16).Value = 1154
$ws.Cells.Item(66, 4).Value = 44550
$ws.Cells.Item(66, 10).Value = 790
$ws.Cells.Item(66, 11).Value = 11000
$ws.Cells.Item(66, 12).Value = 12000
$ws.Cells.Item(66, 13).Value = 11494
$ws.Cells.Item(66, 16).Value = 884
$ws.Cells.Item(67, 4).Value = 44550
$ws.Cells.Item(67, 11).Value = 9000
$ws.Cells.Item(67, 12).Value = 10000
$ws.Cells.Item(67, 13).Value = 9500
$ws.Cells.Item(67, 16).Value = 731
$ws.Cells.Item(68, 4).Value = 44424
$ws.Cells.Item(68, 10).Value = 700
$ws.Cells.Item(68, 11).Value = 13000
$ws.Cells.Item(68, 12).Value = 14000
$ws.Cells.Item(68, 13).Value = 13500
$ws.Cells.Item(68, 16).Value = 1038
$ws.Cells.Item(69, 4).Value = 44424
$ws.Cells.Item(69, 9).Value = "Segunda"
$ws.Cells.Item(69, 10).Value = 430
$ws.Cells.Item(69, 11).Value = 12000
$ws.Cells.Item(69, 12).Value = 12000
$ws.Cells.Item(69, 13).Value = 12000
$ws.Cells.Item(69, 16).Value = 923
$ws.Cells.Item(70, 4).Value = 44396
$ws.Cells.Item(70, 10).Value = 770
$ws.Cells.Item(70, 11).Value = 17000
$ws.Cells.Item(70, 12).Value = 18000
$ws.Cells.Item(70, 13).Value = 17494
$ws.Cells.Item(70, 16).Value = 1346
$ws.Cells.Item(71, 4).Value = 44396
$ws.Cells.Item(71, 11).Value = 16000
$ws.Cells.Item(71, 12).Value = 16000
$ws.Cells.Item(71, 13).Value = 16000
$ws.Cells.Item(71, 16).Value = 1231
$ws.Cells.Item(72, 4).Value = 44235
$ws.Cells.Item(72, 10).Value = 250
$ws.Cells.Item(72, 11).Value = 42000
$ws.Cells.Item(72, 12).Value = 43000
$ws.Cells.Item(72, 13).Value = 42400
$ws.Cells.Item(72, 16).Value = 3262
$ws.Cells.Item(73, 4).Value = 44473
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 1060
$ws.Cells.Item(73, 11).Value = 14000
$ws.Cells.Item(73, 12).Value = 15000
$ws.Cells.Item(73, 13).Value = 14500
$ws.Cells.Item(73, 16).Value = 1115
$ws.Cells.Item(74, 4).Value = 44473
$ws.Cells.Item(74, 9).Value = "Segunda"
$ws.Cells.Item(74, 10).Value = 430
$ws.Cells.Item(74, 11).Value = 11000
$ws.Cells.Item(74, 12).Value = 12000
$ws.Cells.Item(74, 13).Value = 11500
$ws.Cells.Item(74, 16).Value = 885
$ws.Cells.Item(75, 4).Value = 44389
$ws.Cells.Item(75, 10).Value = 700
$ws.Cells.Item(75, 11).Value = 19000
$ws.Cells.Item(75, 12).Value = 20000
$ws.Cells.Item(75, 13).Value = 19500
$ws.Cells.Item(75, 16).Value = 1500
$ws.Cells.Item(76, 4).Value = 44389
$ws.Cells.Item(76, 10).Value = 340
$ws.Cells.Item(76, 11).Value = 17000
$ws.Cells.Item(76, 12).Value = 17000
$ws.Cells.Item(76, 13).Value = 17000
$ws.Cells.Item(76, 16).Value = 1308
$ws.Cells.Item(77, 4).Value = 44249
$ws.Cells.Item(77, 10).Value = 250
$ws.Cells.Item(77, 11).Value = 39000
$ws.Cells.Item(77, 12).Value = 42000
$ws.Cells.Item(77, 13).Value = 40500
$ws.Cells.Item(77, 16).Value = 3115
$ws.Cells.Item(78, 4).Value = 44431
$ws.Cells.Item(79, 4).Value = 44431
$ws.Cells.Item(80, 4).Value = 44382
$ws.Cells.Item(80, 11).Value = 14000
$ws.Cells.Item(80, 12).Value = 15000
$ws.Cells.Item(80, 13).Value = 14500
$ws.Cells.Item(80, 16).Value = 1115
$ws.Cells.Item(81, 4).Value = 44382
$ws.Cells.Item(81, 11).Value = 12000
$ws.Cells.Item(81, 12).Value = 12000
$ws.Cells.Item(81, 13).Value = 12000
$ws.Cells.Item(81, 16).Value = 923
$ws.Cells.Item(82, 4).Value = 44634
$ws.Cells.Item(82, 10).Value = 520
$ws.Cells.Item(82, 11).Value = 16000
$ws.Cells.Item(82, 12).Value = 17000
$ws.Cells.Item(82, 13).Value = 16500
$ws.Cells.Item(82, 16).Value = 1269
$ws.Cells.Item(83, 4).Value = 44487
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 1150
$ws.Cells.Item(83, 11).Value = 14000
$ws.Cells.Item(83, 12).Value = 15000
$ws.Cells.Item(83, 13).Value = 14500
$ws.Cells.Item(83, 16).Value = 1115
$ws.Cells.Item(84, 4).Value = 44487
$ws.Cells.Item(84, 9).Value = "Segunda"
$ws.Cells.Item(84, 10).Value = 610
$ws.Cells.Item(84, 11).Value = 12000
$ws.Cells.Item(84, 12).Value = 12000
$ws.Cells.Item(84, 13).Value = 12000
$ws.Cells.Item(84, 16).Value = 923
$ws.Cells.Item(85, 4).Value = 44578
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 790
$ws.Cells.Item(85, 11).Value = 11000
$ws.Cells.Item(85, 12).Value = 12000
$ws.Cells.Item(85, 13).Value = 11494
$ws.Cells.Item(85, 16).Value = 884
$ws.Cells.Item(86, 4).Value = 44578
$ws.Cells.Item(86, 9).Value = "Segunda"
$ws.Cells.Item(86, 10).Value = 340
$ws.Cells.Item(86, 11).Value = 10000
$ws.Cells.Item(86, 12).Value = 10000
$ws.Cells.Item(86, 13).Value = 10000
$ws.Cells.Item(86, 16).Value = 769
$ws.Cells.Item(87, 4).Value = 44669
$ws.Cells.Item(87, 10).Value = 610
$ws.Cells.Item(87, 11).Value = 10000
$ws.Cells.Item(87, 12).Value = 11000
$ws.Cells.Item(87, 13).Value = 10500
$ws.Cells.Item(87, 16).Value = 808
$ws.Cells.Item(88, 4).Value = 44529
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 790
$ws.Cells.Item(88, 11).Value = 16000
$ws.Cells.Item(88, 12).Value = 18000
$ws.Cells.Item(88, 13).Value = 16987
$ws.Cells.Item(88, 16).Value = 1307
$ws.Cells.Item(89, 4).Value = 44529
$ws.Cells.Item(89, 9).Value = "Segunda"
$ws.Cells.Item(89, 10).Value = 430
$ws.Cells.Item(89, 11).Value = 13000
$ws.Cells.Item(89, 12).Value = 14000
$ws.Cells.Item(89, 13).Value = 13500
$ws.Cells.Item(89, 16).Value = 1038
$ws.Cells.Item(90, 4).Value = 44284
$ws.Cells.Item(90, 10).Value = 340
$ws.Cells.Item(90, 11).Value = 28000
$ws.Cells.Item(90, 12).Value = 30000
$ws.Cells.Item(90, 13).Value = 29000
$ws.Cells.Item(90, 16).Value = 2231
$ws.Cells.Item(91, 4).Value = 44403
$ws.Cells.Item(91, 10).Value = 700
$ws.Cells.Item(91, 11).Value = 16000
$ws.Cells.Item(91, 12).Value = 17000
$ws.Cells.Item(91, 13).Value = 16500
$ws.Cells.Item(91, 16).Value = 1269
$ws.Cells.Item(92, 4).Value = 44403
$ws.Cells.Item(92, 10).Value = 430
$ws.Cells.Item(92, 11).Value = 15000
$ws.Cells.Item(92, 12).Value = 15000
$ws.Cells.Item(92, 13).Value = 15000
$ws.Cells.Item(92, 16).Value = 1154
$ws.Cells.Item(93, 4).Value = 44557
$ws.Cells.Item(93, 10).Value = 970
$ws.Cells.Item(93, 11).Value = 17000
$ws.Cells.Item(93, 12).Value = 18000
$ws.Cells.Item(93, 13).Value = 17495
$ws.Cells.Item(93, 16).Value = 1346
$ws.Cells.Item(94, 4).Value = 44557
$ws.Cells.Item(94, 10).Value = 430
$ws.Cells.Item(94, 11).Value = 16000
$ws.Cells.Item(94, 12).Value = 16000
$ws.Cells.Item(94, 13).Value = 16000
$ws.Cells.Item(94, 16).Value = 1231
$ws.Cells.Item(95, 4).Value = 44305
$ws.Cells.Item(95, 10).Value = 340
$ws.Cells.Item(95, 11).Value = 24000
$ws.Cells.Item(95, 12).Value = 24000
$ws.Cells.Item(95, 13).Value = 24000
$ws.Cells.Item(95, 16).Value = 1846
$ws.Cells.Item(96, 4).Value = 44305
$ws.Cells.Item(96, 10).Value = 160
$ws.Cells.Item(96, 11).Value = 20000
$ws.Cells.Item(96, 12).Value = 20000
$ws.Cells.Item(96, 13).Value = 20000
$ws.Cells.Item(96, 16).Value = 1538
$ws.Cells.Item(97, 4).Value = 44200
$ws.Cells.Item(97, 10).Value = 520
$ws.Cells.Item(97, 11).Value = 30000
$ws.Cells.Item(97, 12).Value = 30000
$ws.Cells.Item(97, 13).Value = 30000
$ws.Cells.Item(97, 16).Value = 2308
$ws.Cells.Item(98, 4).Value = 44200
$ws.Cells.Item(98, 11).Value = 25000
$ws.Cells.Item(98, 12).Value = 25000
$ws.Cells.Item(98, 13).Value = 25000
$ws.Cells.Item(98, 16).Value = 1923
